$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.778.77'
$ws.Range('E2').Value = '  +2.22%  '

$ws.Range('D3').Value = '3.950.63'
$ws.Range('E3').Value = '  +0.88%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '529.21'
$ws.Range('E5').Value = '  +8.10%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.84'
$ws.Range('E6').Value = '  -0.33%  '

$ws.Range('E7').Value = '  -0.25%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('E9').Value = '  -0.38%  '

$ws.Range('E10').Value = '  +4.39%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000342'
$ws.Range('E11').Value = '  -1.11%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.18'
$ws.Range('E12').Value = '  +0.26%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '4.578.52'
$ws.Range('E13').Value = '  +0.88%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.42'
$ws.Range('E14').Value = '  -4.13%  '

$ws.Range('D15').Value = '3.954.80'
$ws.Range('E15').Value = '  +0.94%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.16'
$ws.Range('E16').Value = '  -0.70%  '

$ws.Range('E17').Value = '  -0.30%  '

$ws.Range('E18').Value = '  +7.21%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.91'
$ws.Range('E19').Value = '  -0.01%  '

$ws.Range('D20').Value = '69.725.60'
$ws.Range('E20').Value = '  +2.00%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '435.46'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.39'
$ws.Range('E22').Value = '  -3.65%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.56'
$ws.Range('E23').Value = '  -3.61%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.36'
$ws.Range('E24').Value = '  +0.21%  '

$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.98'
$ws.Range('E25').Value = '  +9.08%  '

$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.80'
$ws.Range('E26').Value = '  +2.11%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.09'
$ws.Range('E27').Value = '  -3.48%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.85'
$ws.Range('E28').Value = '  -4.29%  '

$ws.Range('E29').Value = '  -1.55%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '701.71'
$ws.Range('E30').Value = '  -3.52%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.29'
$ws.Range('E31').Value = '  -3.70%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.127'
$ws.Range('E32').Value = '  -3.09%  '

$ws.Range('E33').Value = '  -1.51%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '67.88'
$ws.Range('E34').Value = '  +10.94%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.442'
$ws.Range('E35').Value = '  +6.77%  '

$ws.Range('E36').Value = '  +1.93%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.96'
$ws.Range('E37').Value = '  -5.62%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '40.36'
$ws.Range('E38').Value = '  -3.96%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.148'
$ws.Range('E39').Value = '  -0.41%  '

$ws.Range('E40').Value = '  -0.27%  '

$ws.Range('E41').Value = '  -0.03%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0484'
$ws.Range('E42').Value = '  +0.42%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.14'
$ws.Range('E43').Value = '  +7.68%  '

$ws.Range('E44').Value = '  -5.72%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.00'
$ws.Range('E45').Value = '  -5.29%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.142'
$ws.Range('E46').Value = '  +0.41%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.35'
$ws.Range('E47').Value = '  +1.99%  '

$ws.Range('D48').Value = '0.0₆0358'
$ws.Range('E48').Value = '  +3.62%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.02'
$ws.Range('E49').Value = '  +6.93%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.35'
$ws.Range('E50').Value = '  -1.89%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.09'
$ws.Range('E51').Value = '  -2.25%  '
